$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("ChatGPT Analysis ", $true, $false, $false, $false, $false, $true, 1, $false, "Chelsea" + [char]0x2019 + "s Research", 2)

$para1 = $d.Paragraphs.Item(1)
Write-Host "Para1 Start=$($para1.Range.Start) End=$($para1.Range.End) Text=[$($para1.Range.Text)]"

$para2 = $d.Paragraphs.Item(2)
Write-Host "Para2 Start=$($para2.Range.Start) End=$($para2.Range.End) Text=[$($para2.Range.Text)]"

# check what happens to the emoji run in the saved file by inspecting the raw paragraph count and structure
Write-Host "Total paragraphs: $($d.Paragraphs.Count)"
Write-Host "Chars count: $($para1.Range.Characters.Count)"
Write-Host "StoryLength: $($para1.Range.StoryLength)"

$emoji = [char]::ConvertFromUtf32(0x1F602)
Write-Host "Searching for emoji char, length=$($emoji.Length)"
$found2 = $d.Content.Find.Execute($emoji, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found emoji: $found2"
if ($found2) {
    Write-Host "Found range: Start=$($d.Content.Find.Parent.Start) End=$($d.Content.Find.Parent.End)"
}

